$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new column before column F (old F shifts to G)
$ws.Columns("F").Insert() | Out-Null

# Set the new column F width
$ws.Columns("F").ColumnWidth = 31.65

# Fill in the new "Player asset" column content
$ws.Range("F2").Value = "主角资源"
$ws.Range("F3").Value = "string"
$ws.Range("F4").Value = "PlayerAsset"
$ws.Range("F5").Value = "Player"

# Update the active selection to F5
$ws.Range("F5").Select() | Out-Null
